# Final OOPS version with Evaluation
# Sync the "Sub-Regression" sheet's descriptive text (column B) with the
# wording already used on the "Supervised" sheet, fixing the
# "analyze" -> "analyse" spelling along the way, and move the
# active-cell selection from B8 down to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sub-Regression")
$ws.Activate()

$ws.Range("B2").Value  = "regression discrete linear observation analyse series "
$ws.Range("B3").Value  = "predict divided approximate range quantity labels class distribution statistics  analyse regression conitnous"
$ws.Range("B4").Value  = "classification category branch segmentation regression analyse  predict series"
$ws.Range("B5").Value  = "regression estimate group separate separate sets"
$ws.Range("B6").Value  = "regression analyse regression predict series forecasting estimate count"
$ws.Range("B7").Value  = "classification classification classification category category category separate decision segmentation segregate"
$ws.Range("B8").Value  = "regression analyse predict segmentation classification classification classification classification classification classification classification classification classification "
$ws.Range("B9").Value  = "regression regression analyse predict segmentation separate"
$ws.Range("B10").Value = "regression analyse regression category category category category separate segmentation"

$ws.Range("B11").Select()
